# Generate Report for Handoff
#
# Moves the localization job from "In Translation" to "Ready for handoff"
# and refreshes the handoff timestamps that the status report tracks:
#   - Overview sheet:  zh-cn / de-de status cells + "Latest HO Xliff
#                       Generate Date"
#   - zh-cn sheet:      Status + "Latest Handoff Datetime"
#   - de-de sheet:      Status + "Latest Handoff Datetime"
#
# Column C/E/F also grow a bit because "Ready for handoff" is longer than
# "In Translation", so the report's column widths are refreshed too.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = "2016-08-29 21:13:19"

# --- zh-cn sheet ------------------------------------------------------
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = "2016-08-29 21:13:14"

# --- de-de sheet ------------------------------------------------------
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = "2016-08-29 21:13:19"

# --- Widen the Status columns to fit the new, longer text -----------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # zh-cn status column
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # de-de status column
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33        # Status column
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33        # Status column
